$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Inventory
$ws.Range("C4").Value = 281000000.0
$ws.Range("D4").Value = 245000000.0
$ws.Range("E4").Value = 219000000.0
$ws.Range("F4").Value = 187000000.0
$ws.Range("G4").Value = 168000000.0

# Row 6 - Total current assets
$ws.Range("B6").Value = 8539307000.0

# Row 12 - Long-term assets (Other)
$ws.Range("B12").Value = 49262000.0

# Row 13 - Total non-current assets
$ws.Range("B13").Value = 3575752000.0

# Row 15 - Accounts Payable
$ws.Range("B15").Value = 127839000.0
$ws.Range("C15").Value = 155000000.0
$ws.Range("D15").Value = 108000000.0
$ws.Range("E15").Value = 101000000.0
$ws.Range("F15").Value = 101000000.0
$ws.Range("G15").Value = 88000000.0

# Row 16 - Accrued Expenses
$ws.Range("B16").Value = 1532037000.0

# Row 21 - Other current liabilities
$ws.Range("B21").Value = 284174000.0

# Row 22 - Total current liabilities
$ws.Range("B22").Value = 1944050000.0

# Row 23 - Long Term Debt (Total)
$ws.Range("B23").Value = 530330000.0

# Row 25 - Long Term Tax Liability (Deferred)
$ws.Range("C25").Value = -883000000.0
$ws.Range("D25").Value = -1148000000.0
$ws.Range("E25").Value = -1215000000.0
$ws.Range("F25").Value = -1148000000.0
$ws.Range("G25").Value = -1191000000.0

# Row 27 - Non-current Liabilities (Other)
$ws.Range("B27").Value = 291958000.0

# Row 28 - Total non-current liabilities
$ws.Range("B28").Value = 1190755000.0

# Row 30 - Additional Paid In Capital
$ws.Range("B30").Value = 7499161000.0

# Row 31 - Common Stock (Net)
$ws.Range("B31").Value = 2588000.0

# Row 32 - Retained Earnings
$ws.Range("B32").Value = 1511807000.0

# Row 33 - Common Equity (Total)
$ws.Range("B33").Value = 8980254000.0

# Row 38 - Net Debt
$ws.Range("B38").Value = -6923968000.0
